$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the "index" column header to "i" (this also renames the table column
# since A1 is the header cell of the "testdata" table's first column).
$ws.Range("A1").Value = "i"

# Shrink column A's width now that the header text is shorter.
$ws.Columns.Item(1).ColumnWidth = 4

# Decrement every data value in column A by 1 (old 1,2,3,... -> new 0,1,2,...).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row()
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $v = $cell.Value()
    $cell.Value = $v - 1
}
